$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.390.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.002.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.57%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.001.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("E13").Value = "  -4.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.122"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.495.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.405.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.004.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.81%  "
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.59%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("E34").Value = "  -5.50%  "
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0789"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "407.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.63%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.112"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.277"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.771.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0351"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.04%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.34%  "
